$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 2271
    $ws.Range("F3").Value = 1719
    $ws.Range("F6").Value = 827
    $ws.Range("F8").Value = 5839
}
